$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.457.15"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.872.29"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D5").Value = "'243.79"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'0.7053"
$ws.Range("E6").Value = "  -2.33%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.07946"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").Value = "'0.3144"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'24.56"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("D11").Value = "'0.07805"
$ws.Range("E11").Value = "  -4.44%  "
$ws.Range("D12").Value = "1.873.17"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "'93.76"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "'5.169"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "'0.7034"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "'6.518"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "'0.000008562"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").Value = "29.481.76"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'252.98"
$ws.Range("E19").Value = "  +3.49%  "
$ws.Range("D20").Value = "2.142.91"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("E21").Value = "  -1.65%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").Value = "'7.627"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "'0.1542"
$ws.Range("E25").Value = "  -3.94%  "
$ws.Range("D26").Value = "'9.006"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").Value = "'161.25"
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("D28").Value = "'18.78"
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("D29").Value = "'1.542"
$ws.Range("E29").Value = "  +2.41%  "
$ws.Range("D30").Value = "'4.309"
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("D31").Value = "'4.264"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  -2.85%  "
$ws.Range("D33").Value = "'0.05287"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").Value = "'1.895"
$ws.Range("E34").Value = "  -2.35%  "
$ws.Range("D35").Value = "'0.7606"
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("D36").Value = "'1.184"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "'2.706"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "'0.01878"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").Value = "1.277.28"
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("D40").Value = "'2.760"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "'0.9002"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").Value = "'109.89"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").Value = "'5.990"
$ws.Range("E43").Value = "  -7.05%  "
$ws.Range("D44").Value = "'71.01"
$ws.Range("E44").Value = "  -4.67%  "
$ws.Range("D45").Value = "'0.9999"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").Value = "2.044.13"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("E47").Value = "  -3.17%  "
$ws.Range("D48").Value = "'9.660"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("D50").Value = "'0.5173"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("E51").Value = "  -0.78%  "
